# Workbook: processor_wla_validation.xlsx
# Commit message: "Added best fit scripts for merom and penryn test cases"
#
# This script:
#  1. Renames "Merom2" to "65nm Merom"
#  2. Adds a "best fit" data block (rows 24-33) to the "65nm Merom" sheet,
#     mirroring the block already present on the "45nm Penryn" sheet.
#  3. Reorders the sheet tabs so "45nm Penryn" and "32nm Sandy Bridge (standard)"
#     move up next to the other sheets that were being worked on.
#  4. Restores/updates view state (selections, active tab) left over from the edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "Merom2" -> "65nm Merom"
# ---------------------------------------------------------------------------
$merom = $wb.Worksheets.Item("Merom2")
$merom.Name = "65nm Merom"

# ---------------------------------------------------------------------------
# 2. Add the new "best fit" block to the Merom sheet (rows 24-33)
# ---------------------------------------------------------------------------
$merom.Range("A24").Value = "Metal"
$merom.Range("B24").Value = "Pitch (nm)"
$merom.Range("C24").Value = "Thickness (nm)"
$merom.Range("D24").Value = "AR"
$merom.Range("E24").Value = "Width"
$merom.Range("F24").Value = "width_frac"
$merom.Range("A24:D24").Font.Bold = $true
$merom.Range("A24:D24").HorizontalAlignment = -4108
$merom.Range("E24:F24").Font.Bold = $true

$data = @(
    @(1, 210, 170, 1.6),
    @(2, 210, 190, 1.8),
    @(3, 220, 200, 1.8),
    @(4, 280, 250, 1.8),
    @(5, 330, 300, 1.8),
    @(6, 480, 430, 1.8),
    @(7, 720, 650, 1.8),
    @(8, 1080, 975, 1.8)
)

$row = 25
foreach ($entry in $data) {
    $merom.Range("A$row").Value = $entry[0]
    $merom.Range("B$row").Value = $entry[1]
    $merom.Range("C$row").Value = $entry[2]
    $merom.Range("D$row").Value = $entry[3]
    $merom.Range("A$row").HorizontalAlignment = -4108
    $merom.Range("C$row").HorizontalAlignment = -4108
    $merom.Range("D$row").HorizontalAlignment = -4108
    $row = $row + 1
}

$merom.Range("E25").Formula = "=C25/D25"
$merom.Range("F25").Formula = "=E25/B25"
$merom.Range("E26:E32").Formula = "=C26/D26"
$merom.Range("F26:F32").Formula = "=E26/B26"

# Trailing (mostly empty) row 33 - only A/C/D carry the centered style
$merom.Range("A33").HorizontalAlignment = -4108
$merom.Range("C33").HorizontalAlignment = -4108
$merom.Range("D33").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 3. Reorder sheet tabs
#    Target order:
#    Merom 65nm Data, 65nm Merom, 45nm Penryn, Sandy Bridge EP-4 32nm Data,
#    32nm Sandy Bridge (standard), 22nm - Ivy Bridge EP10,
#    22nm Ivy Bridge (standard), Sheet1
# ---------------------------------------------------------------------------
$penryn = $wb.Worksheets.Item("45nm Penryn")
$penryn.Activate()
$penryn.Range("A7:F16").Select()

$sandyStd = $wb.Worksheets.Item("32nm Sandy Bridge (standard)")

$penryn.Move($wb.Worksheets.Item("Sandy Bridge EP-4 32nm Data"))
$sandyStd.Move($wb.Worksheets.Item("22nm - Ivy Bridge EP10"))

# ---------------------------------------------------------------------------
# 4. Restore view state on the Merom sheet and make it the active tab
# ---------------------------------------------------------------------------
$merom.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$merom.Range("D33").Select()

Write-Host "done"
